$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text-formatted so numeric-looking strings
# (e.g. "24.00", "0.540") are not silently coerced into numbers and lose
# their trailing zeros / multi-dot grouping, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.492.65"
$ws.Range("E2").Value = "  +4.79%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.600.44"
$ws.Range("E3").Value = "  +2.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "214.97"
$ws.Range("E5").Value = "  +2.13%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.498"
$ws.Range("E6").Value = "  +1.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Solana
$ws.Range("D8").Value = "24.01"
$ws.Range("E8").Value = "  +8.70%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.252"
$ws.Range("E9").Value = "  +1.25%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").Value = "  +0.90%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +2.32%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.829.38"
$ws.Range("E12").Value = "  +2.58%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.592.75"
$ws.Range("E13").Value = "  +2.00%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.534"
$ws.Range("E14").Value = "  +3.43%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "3.77"
$ws.Range("E15").Value = "  +0.50%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "28.525.78"
$ws.Range("E16").Value = "  +5.04%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.31"
$ws.Range("E17").Value = "  +2.51%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "231.23"
$ws.Range("E18").Value = "  +6.96%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +1.27%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +1.75%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.02%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.39"
$ws.Range("E23").Value = "  +2.04%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.34%  "

# Row 25 - Monero
$ws.Range("D25").Value = "152.35"
$ws.Range("E25").Value = "  -0.05%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "15.26"
$ws.Range("E26").Value = "  +1.62%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "6.63"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "0.107"
$ws.Range("E28").Value = "  +0.72%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.00%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.45%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.0475"
$ws.Range("E31").Value = "  +1.43%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  +0.23%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  +0.06%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.399.84"
$ws.Range("E34").Value = "  -3.14%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -0.53%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  -4.44%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  +0.96%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  +1.35%  "

# Row 39 - swapped: now MXToken
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.52"
$ws.Range("E39").Value = "  +8.37%  "

# Row 40 - swapped: now ImmutableX
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.540"
$ws.Range("E40").Value = "  -2.46%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Value = "0.820"
$ws.Range("E41").Value = "  +1.50%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "5.72"
$ws.Range("E42").Value = "  -2.94%  "

# Row 43 - PaxDollar
$ws.Range("E43").Value = "  -0.01%  "

# Row 44 - WEMIXToken
$ws.Range("D44").Value = "0.985"
$ws.Range("E44").Value = "  -0.82%  "

# Row 45 - RenderToken
$ws.Range("E45").Value = "  +6.35%  "

# Row 46 - Aave
$ws.Range("D46").Value = "64.71"
$ws.Range("E46").Value = "  +1.06%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.737.93"
$ws.Range("E47").Value = "  +2.51%  "

# Row 48 - mCoin
$ws.Range("E48").Value = "  +0.58%  "

# Row 49 - Quant
$ws.Range("D49").Value = "87.45"
$ws.Range("E49").Value = "  +2.24%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +5.95%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.0529"
$ws.Range("E51").Value = "  +0.96%  "
